$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3290
$ws.Range("I18").Value = 2490
$ws.Range("J18").Value = 4490
$ws.Range("K18").Value = 2490
$ws.Range("L18").Value = 4490
$ws.Range("M18").Value = -2206
$ws.Range("N18").Value = -5058
$ws.Range("H40").Value = 83339656
$ws.Range("I40").Value = 1489
$ws.Range("J40").Value = 125008744
$ws.Range("K40").Value = 1489
$ws.Range("L40").Value = 125008744
$ws.Range("M40").Value = -1314
$ws.Range("N40").Value = -125009094
$ws.Range("H103").Value = 50001730
$ws.Range("I103").Value = 800
$ws.Range("J103").Value = 62501960
$ws.Range("K103").Value = 2400
$ws.Range("L103").Value = 187505880
$ws.Range("M103").Value = -1814
$ws.Range("N103").Value = -187507052
$ws.Range("H111").Value = 12516
$ws.Range("J111").Value = 12516
$ws.Range("L111").Value = 37548
$ws.Range("N111").Value = -43682
$ws.Range("H116").Value = 14093.895
$ws.Range("I116").Value = 5214.5
$ws.Range("J116").Value = 18192.076
$ws.Range("K116").Value = 5214.5
$ws.Range("L116").Value = 18192.076
$ws.Range("M116").Value = -1772.5
$ws.Range("N116").Value = -25076.076
$ws.Range("H132").Value = 1678.6765
$ws.Range("I132").Value = 1718.9395
$ws.Range("J132").Value = 350
$ws.Range("K132").Value = 5156.818499999999
$ws.Range("L132").Value = 1050
$ws.Range("M132").Value = -2626.818499999999
$ws.Range("N132").Value = -6110
$ws.Range("H137").Value = 3035.2632
$ws.Range("I137").Value = 1960
$ws.Range("J137").Value = 3817.2727
$ws.Range("K137").Value = 5880
$ws.Range("L137").Value = 11451.8181
$ws.Range("M137").Value = -3330
$ws.Range("N137").Value = -16551.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 52635372
$ws.Range("I4").Value = 4007.1428
$ws.Range("J4").Value = 200003200
$ws.Range("K4").Value = 4007.1428
$ws.Range("L4").Value = 200003200
$ws.Range("M4").Value = -3891.1428
$ws.Range("N4").Value = -200003432
$ws.Range("H32").Value = 6346.527
$ws.Range("I32").Value = 5901.2695
$ws.Range("J32").Value = 16253.5
$ws.Range("K32").Value = 5901.2695
$ws.Range("L32").Value = 16253.5
$ws.Range("M32").Value = -5614.2695
$ws.Range("N32").Value = -16827.5
$ws.Range("H45").Value = 3822.6667
$ws.Range("I45").Value = 2586
$ws.Range("K45").Value = 2586
$ws.Range("M45").Value = -2209

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 7262
$ws.Range("I22").Value = 1272.2
$ws.Range("K22").Value = 1272.2
$ws.Range("M22").Value = -1099.2
$ws.Range("H134").Value = 11112268
$ws.Range("I134").Value = 815.7143
$ws.Range("K134").Value = 2447.1429
$ws.Range("M134").Value = 87.85710000000017

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7143963
$ws.Range("I16").Value = 11111943
$ws.Range("J16").Value = 1599.8
$ws.Range("K16").Value = 11111943
$ws.Range("L16").Value = 1599.8
$ws.Range("M16").Value = -11111656
$ws.Range("N16").Value = -2173.8
$ws.Range("H31").Value = 76927840
$ws.Range("I31").Value = 166670540
$ws.Range("J31").Value = 5525
$ws.Range("K31").Value = 166670540
$ws.Range("L31").Value = 5525
$ws.Range("M31").Value = -166670245
$ws.Range("N31").Value = -6115
$ws.Range("H34").Value = 76927840
$ws.Range("I34").Value = 166670540
$ws.Range("J34").Value = 5525
$ws.Range("K34").Value = 166670540
$ws.Range("L34").Value = 5525
$ws.Range("M34").Value = -166670338
$ws.Range("N34").Value = -5929
$ws.Range("H58").Value = 1737.6471
$ws.Range("I58").Value = 1084.6072
$ws.Range("J58").Value = 4785.1665
$ws.Range("K58").Value = 1084.6072
$ws.Range("L58").Value = 4785.1665
$ws.Range("M58").Value = -881.6071999999999
$ws.Range("N58").Value = -5191.1665
$ws.Range("H107").Value = 2929.7778
$ws.Range("I107").Value = 796
$ws.Range("K107").Value = 796
$ws.Range("M107").Value = 1124
$ws.Range("H113").Value = 7143963
$ws.Range("I113").Value = 11111943
$ws.Range("J113").Value = 1599.8
$ws.Range("K113").Value = 11111943
$ws.Range("L113").Value = 1599.8
$ws.Range("M113").Value = -11109773
$ws.Range("N113").Value = -5939.8
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
$ws.Range("H134").Value = 4140.25
$ws.Range("I134").Value = 3266.5
$ws.Range("J134").Value = 5014
$ws.Range("K134").Value = 9799.5
$ws.Range("L134").Value = 15042
$ws.Range("M134").Value = -7264.5
$ws.Range("N134").Value = -20112
$ws.Range("H136").Value = 1737.6471
$ws.Range("I136").Value = 1084.6072
$ws.Range("J136").Value = 4785.1665
$ws.Range("K136").Value = 3253.8216
$ws.Range("L136").Value = 14355.4995
$ws.Range("M136").Value = -703.8215999999998
$ws.Range("N136").Value = -19455.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 3911.2
$ws.Range("I26").Value = 799.8570999999999
$ws.Range("K26").Value = 2399.5713
$ws.Range("M26").Value = -2111.5713
$ws.Range("H38").Value = 39.666668
$ws.Range("I38").Value = 32.3
$ws.Range("J38").Value = 76.5
$ws.Range("K38").Value = 96.89999999999999
$ws.Range("L38").Value = 229.5
$ws.Range("M38").Value = 250.1
$ws.Range("N38").Value = -923.5
$ws.Range("H114").Value = 4823.1816
$ws.Range("I114").Value = 113.28571
$ws.Range("J114").Value = 13065.5
$ws.Range("K114").Value = 339.85713
$ws.Range("L114").Value = 39196.5
$ws.Range("M114").Value = 2914.14287
$ws.Range("N114").Value = -45704.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14594.6
$ws.Range("I70").Value = 14937.75
$ws.Range("K70").Value = 14937.75
$ws.Range("M70").Value = -14667.75
$ws.Range("H73").Value = 14594.6
$ws.Range("I73").Value = 14937.75
$ws.Range("K73").Value = 14937.75
$ws.Range("M73").Value = -14001.75
$ws.Range("H102").Value = 4133.875
$ws.Range("I102").Value = 4215.25
$ws.Range("J102").Value = 4052.5
$ws.Range("K102").Value = 4215.25
$ws.Range("L102").Value = 4052.5
$ws.Range("M102").Value = -2593.25
$ws.Range("N102").Value = -7296.5
$ws.Range("H113").Value = 1547372.1
$ws.Range("I113").Value = 2170
$ws.Range("K113").Value = 2170
$ws.Range("M113").Value = 0
$ws.Range("H126").Value = 5664.4443
$ws.Range("I126").Value = 5176.1
$ws.Range("K126").Value = 15528.3
$ws.Range("M126").Value = -13058.3
$ws.Range("H132").Value = 11116769
$ws.Range("I132").Value = 6564
$ws.Range("J132").Value = 33337180
$ws.Range("K132").Value = 19692
$ws.Range("L132").Value = 100011540
$ws.Range("M132").Value = -17162
$ws.Range("N132").Value = -100016600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3043.5527
$ws.Range("I61").Value = 2718.6428
$ws.Range("K61").Value = 2718.6428
$ws.Range("M61").Value = -2516.6428
$ws.Range("H113").Value = 3043.5527
$ws.Range("I113").Value = 2718.6428
$ws.Range("K113").Value = 2718.6428
$ws.Range("M113").Value = -548.6428000000001
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H132").Value = 3784.1177
$ws.Range("I132").Value = 2022.2727
$ws.Range("J132").Value = 7014.1665
$ws.Range("K132").Value = 6066.8181
$ws.Range("L132").Value = 21042.4995
$ws.Range("M132").Value = -3536.8181
$ws.Range("N132").Value = -26102.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5002498.5
$ws.Range("I132").Value = 4995
$ws.Range("J132").Value = 6668333
$ws.Range("K132").Value = 14985
$ws.Range("L132").Value = 20004999
$ws.Range("M132").Value = -12455
$ws.Range("N132").Value = -20010059
$ws.Range("H136").Value = 183395.48
$ws.Range("I136").Value = 1524.9535
$ws.Range("J136").Value = 835098.25
$ws.Range("K136").Value = 4574.860500000001
$ws.Range("L136").Value = 2505294.75
$ws.Range("M136").Value = -2024.860500000001
$ws.Range("N136").Value = -2510394.75
